# Updated legacy GSC export data: the oldest/stale row (2025-09-23, a
# duplicate leading row with blank indexed counts) has been removed from
# the "Chart" sheet's data table, shifting every subsequent day up by one
# row (A1:D88 -> A1:D87). The "Table" and "Metadata" sheets are untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$ws.Rows.Item(2).Delete()
